# Updated cryptos list -- price (column D) and 1h volume % (column E) refresh.
# Values that look like plain decimal numbers are prefixed with a leading
# apostrophe so Excel stores them as text (matching the sheet's existing
# inlineStr/text cell type for these columns) instead of silently coercing
# them to a Number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.316.38"
$ws.Range("E2").Value = "  -1.43%  "

$ws.Range("D3").Value = "2.036.22"
$ws.Range("E3").Value = "  -0.31%  "

$ws.Range("E4").Value = "  -0.05%  "

$ws.Range("D5").Value = "'244.45"
$ws.Range("E5").Value = "  -0.08%  "

$ws.Range("D6").Value = "'0.654"
$ws.Range("E6").Value = "  +0.36%  "

$ws.Range("E7").Value = "  +0.04%  "

$ws.Range("D8").Value = "'52.61"
$ws.Range("E8").Value = "  -7.43%  "

$ws.Range("D9").Value = "'60.86"
$ws.Range("E9").Value = "  +4.07%  "

$ws.Range("D10").Value = "'0.357"
$ws.Range("E10").Value = "  -2.58%  "

$ws.Range("D11").Value = "'0.0735"
$ws.Range("E11").Value = "  -4.64%  "

$ws.Range("E12").Value = "  -3.96%  "

$ws.Range("D13").Value = "'0.921"
$ws.Range("E13").Value = "  +6.58%  "

$ws.Range("D14").Value = "'14.32"
$ws.Range("E14").Value = "  -4.74%  "

$ws.Range("D15").Value = "2.334.49"
$ws.Range("E15").Value = "  -0.41%  "

$ws.Range("D16").Value = "'5.30"
$ws.Range("E16").Value = "  -4.89%  "

$ws.Range("D17").Value = "2.029.01"
$ws.Range("E17").Value = "  -2.10%  "

$ws.Range("D18").Value = "36.264.33"
$ws.Range("E18").Value = "  -1.36%  "

$ws.Range("D19").Value = "'16.74"
$ws.Range("E19").Value = "  -5.96%  "

$ws.Range("D20").Value = "'70.75"
$ws.Range("E20").Value = "  -3.12%  "

$ws.Range("D21").Value = "0.0₃0839"
$ws.Range("E21").Value = "  -4.69%  "

$ws.Range("D22").Value = "'234.86"
$ws.Range("E22").Value = "  -0.09%  "

$ws.Range("E23").Value = "  -4.47%  "

$ws.Range("E24").Value = "  -0.08%  "

$ws.Range("E25").Value = "  -3.13%  "

$ws.Range("E26").Value = "  +2.02%  "

$ws.Range("D27").Value = "'163.33"
$ws.Range("E27").Value = "  -2.81%  "

$ws.Range("D28").Value = "'9.01"
$ws.Range("E28").Value = "  -11.29%  "

$ws.Range("D29").Value = "'19.60"
$ws.Range("E29").Value = "  -0.99%  "

$ws.Range("D30").Value = "'0.119"
$ws.Range("E30").Value = "  -2.79%  "

$ws.Range("E31").Value = "  +6.73%  "

$ws.Range("E32").Value = "  -9.97%  "

$ws.Range("E33").Value = "  -3.72%  "

$ws.Range("D34").Value = "'4.33"
$ws.Range("E34").Value = "  -6.96%  "

$ws.Range("E35").Value = "  -0.01%  "

$ws.Range("D36").Value = "'0.0853"
$ws.Range("E36").Value = "  +4.73%  "

$ws.Range("D37").Value = "'1.82"
$ws.Range("E37").Value = "  -0.39%  "

$ws.Range("E38").Value = "  -5.19%  "

$ws.Range("D39").Value = "'4.89"
$ws.Range("E39").Value = "  -4.01%  "

$ws.Range("E40").Value = "  -7.27%  "

$ws.Range("D41").Value = "'2.88"
$ws.Range("E41").Value = "  -5.04%  "

$ws.Range("D42").Value = "'0.0210"
$ws.Range("E42").Value = "  -5.13%  "

$ws.Range("E43").Value = "  -4.29%  "

$ws.Range("D44").Value = "'91.78"
$ws.Range("E44").Value = "  -4.25%  "

$ws.Range("D45").Value = "'0.0884"
$ws.Range("E45").Value = "  -5.23%  "

$ws.Range("D46").Value = "1.367.87"
$ws.Range("E46").Value = "  +5.54%  "

$ws.Range("D47").Value = "'15.45"
$ws.Range("E47").Value = "  -7.55%  "

$ws.Range("E48").Value = "  +9.35%  "

$ws.Range("D49").Value = "'2.91"
$ws.Range("E49").Value = "  +2.09%  "

$ws.Range("D50").Value = "2.220.83"
$ws.Range("E50").Value = "  -0.34%  "

$ws.Range("E51").Value = "  -4.62%  "
